$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 0.1172413793103448
$ws.Range("C2").Value = 0.6275862068965518
$ws.Range("J2").Value = 0.01724137931034483
$ws.Range("P2").Value = 0.1482758620689655
$ws.Range("S2").Value = 0.0896551724137931
$ws.Range("B3").Value = 0.015
$ws.Range("C3").Value = 0.02
$ws.Range("J3").Value = 0.045
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.17
$ws.Range("J4").Value = 0.1395348837209302
$ws.Range("P4").Value = 0.6744186046511628
$ws.Range("S4").Value = 0.186046511627907
$ws.Range("B6").Value = 0.06696428571428571
$ws.Range("D6").Value = 0.01785714285714286
$ws.Range("F6").Value = 0.05357142857142857
$ws.Range("J6").Value = 0.2723214285714285
$ws.Range("O6").Value = 0.01785714285714286
$ws.Range("Q6").Value = 0.15625
$ws.Range("R6").Value = 0.04017857142857143
$ws.Range("S6").Value = 0.375
$ws.Range("B7").Value = 0.1024390243902439
$ws.Range("D7").Value = 0.004878048780487805
$ws.Range("F7").Value = 0.04390243902439024
$ws.Range("J7").Value = 0.07317073170731707
$ws.Range("O7").Value = 0.1024390243902439
$ws.Range("Q7").Value = 0.1707317073170732
$ws.Range("R7").Value = 0.1073170731707317
$ws.Range("S7").Value = 0.3951219512195122
$ws.Range("B8").Value = 0.1282051282051282
$ws.Range("D8").Value = 0.02169625246548323
$ws.Range("F8").Value = 0.0650887573964497
$ws.Range("J8").Value = 0.1222879684418146
$ws.Range("O8").Value = 0.01577909270216963
$ws.Range("Q8").Value = 0.1873767258382643
$ws.Range("R8").Value = 0.07297830374753451
$ws.Range("S8").Value = 0.3865877712031558
$ws.Range("B9").Value = 0.09473684210526316
$ws.Range("D9").Value = 0.02105263157894737
$ws.Range("F9").Value = 0.07368421052631578
$ws.Range("J9").Value = 0.08421052631578947
$ws.Range("O9").Value = 0.02105263157894737
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.05789473684210526
$ws.Range("S9").Value = 0.4473684210526316
$ws.Range("B10").Value = 0.08895478131949593
$ws.Range("D10").Value = 0.01853224610822832
$ws.Range("E10").Value = 0.002223869532987398
$ws.Range("F10").Value = 0.06300963676797627
$ws.Range("J10").Value = 0.1779095626389919
$ws.Range("O10").Value = 0.02446256486286138
$ws.Range("Q10").Value = 0.1979243884358784
$ws.Range("R10").Value = 0.083765752409192
$ws.Range("S10").Value = 0.3432171979243884
$ws.Range("G11").Value = 0.1619718309859155
$ws.Range("J11").Value = 0.09507042253521127
$ws.Range("K11").Value = 0.1549295774647887
$ws.Range("L11").Value = 0.5880281690140845
$ws.Range("G12").Value = 0.727810650887574
$ws.Range("J12").Value = 0.2130177514792899
$ws.Range("K12").Value = 0.005917159763313609
$ws.Range("L12").Value = 0.02958579881656805
$ws.Range("S12").Value = 0.02366863905325444
$ws.Range("F13").Value = 0.01923076923076923
$ws.Range("G13").Value = 0.7115384615384616
$ws.Range("J13").Value = 0.2307692307692308
$ws.Range("S13").Value = 0.03846153846153846
$ws.Range("F15").Value = 0.01834862385321101
$ws.Range("H15").Value = 0.1651376146788991
$ws.Range("I15").Value = 0.05504587155963303
$ws.Range("J15").Value = 0.3440366972477064
$ws.Range("K15").Value = 0.06880733944954129
$ws.Range("M15").Value = 0.01834862385321101
$ws.Range("O15").Value = 0.04587155963302753
$ws.Range("S15").Value = 0.2844036697247707
$ws.Range("F16").Value = 0.08796296296296297
$ws.Range("H16").Value = 0.2083333333333333
$ws.Range("I16").Value = 0.07870370370370371
$ws.Range("J16").Value = 0.3703703703703703
$ws.Range("K16").Value = 0.08333333333333333
$ws.Range("M16").Value = 0.02314814814814815
$ws.Range("O16").Value = 0.02777777777777778
$ws.Range("S16").Value = 0.1203703703703704
$ws.Range("F17").Value = 0.01271186440677966
$ws.Range("H17").Value = 0.173728813559322
$ws.Range("I17").Value = 0.1059322033898305
$ws.Range("J17").Value = 0.4004237288135593
$ws.Range("K17").Value = 0.08898305084745763
$ws.Range("M17").Value = 0.02754237288135593
$ws.Range("O17").Value = 0.05932203389830509
$ws.Range("S17").Value = 0.1313559322033898
$ws.Range("F18").Value = 0.02072538860103627
$ws.Range("H18").Value = 0.1606217616580311
$ws.Range("I18").Value = 0.07253886010362694
$ws.Range("J18").Value = 0.4145077720207254
$ws.Range("K18").Value = 0.1398963730569948
$ws.Range("M18").Value = 0.0310880829015544
$ws.Range("O18").Value = 0.05181347150259067
$ws.Range("S18").Value = 0.1088082901554404
$ws.Range("F19").Value = 0.01160990712074303
$ws.Range("H19").Value = 0.2360681114551084
$ws.Range("I19").Value = 0.07662538699690402
$ws.Range("J19").Value = 0.3676470588235294
$ws.Range("K19").Value = 0.1075851393188854
$ws.Range("M19").Value = 0.02012383900928793
$ws.Range("N19").Value = 0.0007739938080495357
$ws.Range("O19").Value = 0.05882352941176471
$ws.Range("S19").Value = 0.1207430340557276
